$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.60199999999999
$ws.Range("A6").Value = -22.54300000000001
$ws.Range("A7").Value = -20.02129999999999
$ws.Range("C7").Value = -12.0513
$ws.Range("A8").Value = -22.24060000000001
$ws.Range("C11").Value = -11.5425
$ws.Range("C12").Value = -10.57079999999999
$ws.Range("E12").Value = 17.49660000000002
$ws.Range("E13").Value = 16.67420000000001
$ws.Range("E14").Value = 17.01610000000001
$ws.Range("C15").Value = -14.70449999999999
$ws.Range("A16").Value = -22.05740000000001
$ws.Range("E16").Value = 16.157
$ws.Range("E19").Value = 16.35619999999999
$ws.Range("A20").Value = -19.5844
$ws.Range("C20").Value = -12.7674
$ws.Range("E20").Value = 16.25229999999998
$ws.Range("A21").Value = -20.11039999999998
$ws.Range("C21").Value = -11.99200000000001
$ws.Range("C22").Value = -12.07150000000001
$ws.Range("E22").Value = 17.04660000000003
$ws.Range("C23").Value = -12.11330000000001
$ws.Range("A28").Value = -21.83769999999999
$ws.Range("A29").Value = -21.52859999999998
$ws.Range("C29").Value = -11.57630000000001
$ws.Range("A30").Value = -21.6039
$ws.Range("A32").Value = -21.29599999999999
$ws.Range("C34").Value = -11.23580000000002
$ws.Range("E36").Value = 16.07970000000001
$ws.Range("A40").Value = -19.83749999999999
$ws.Range("C42").Value = -12.20480000000001
$ws.Range("C43").Value = -12.97219999999999
$ws.Range("E43").Value = 17.31180000000001
$ws.Range("C44").Value = -13.9007
$ws.Range("C45").Value = -13.95699999999998
$ws.Range("A46").Value = -21.98889999999999
$ws.Range("C46").Value = -12.8199
$ws.Range("E46").Value = 17.00509999999999
$ws.Range("C50").Value = -14.23869999999999
$ws.Range("E50").Value = 16.32889999999999
$ws.Range("A51").Value = -21.72889999999999
$ws.Range("C51").Value = -11.0086
$ws.Range("A52").Value = -22.23850000000001
$ws.Range("A57").Value = -21.9003
$ws.Range("C57").Value = -12.37119999999999
$ws.Range("A59").Value = -22.32130000000001
$ws.Range("A62").Value = -22.19920000000002
$ws.Range("C65").Value = -13.55599999999999
$ws.Range("A66").Value = -21.81619999999999
$ws.Range("C66").Value = -11.1205
$ws.Range("C67").Value = -11.0054
$ws.Range("A73").Value = -20.34229999999999
$ws.Range("A74").Value = -22.10939999999999
$ws.Range("E76").Value = 16.4832
$ws.Range("A77").Value = -20.02319999999999
$ws.Range("C79").Value = -11.82020000000001
$ws.Range("C84").Value = -13.44959999999999
$ws.Range("C87").Value = -13.7257
$ws.Range("A92").Value = -21.76219999999999
$ws.Range("C92").Value = -11.1667
$ws.Range("E95").Value = 18.09030000000002
$ws.Range("C97").Value = -10.8777
$ws.Range("E97").Value = 16.5142
$ws.Range("E99").Value = 16.6104
$ws.Range("A100").Value = -22.2961
